$wb = $excel.ActiveWorkbook

# --- Design sheet: fill in the Design test result for item 31 (row 33) ---
# Previously the comment said "Ergebnis wird morgen nachgetragen" (result to be
# added tomorrow). Now the actual result has come in: Ergebnis Design = 1 (pass)
# and the placeholder comment is cleared.
$design = $wb.Worksheets.Item("Design")
$design.Range("D33").Value = 1
$design.Range("E33").Value = ""
[void]$design.Rows.Item(33).EntireRow.AutoFit()
[void]$design.Range("D34").Select()

# --- Abnahmetest (consolidated) sheet ---
$ws = $wb.Worksheets.Item("Abnahmetest")

# The calculated-column cells for row 33 (item 31) are typed over directly,
# which breaks their table formula for just this row and stores literal
# values instead (Design result = 1, Design comment now empty).
$ws.Range("H33").Value = 1
$ws.Range("I33").Value = ""

# Mark every row's "Ergebnis Überarbeitung" (revision result) column as done.
$ws.Range("J3:J33").Value = 1

[void]$ws.Range("I36").Select()

[void]$ws.Activate()
